$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (B5:AH5) values to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 1.81
$ws.Range("C5").Value = 1.22
$ws.Range("D5").Value = 0.61
$ws.Range("E5").Value = 4.17
$ws.Range("F5").Value = 2.81
$ws.Range("G5").Value = 1.37
$ws.Range("H5").Value = 8.71
$ws.Range("I5").Value = 2.33
$ws.Range("J5").Value = 0.94
$ws.Range("K5").Value = 1.1
$ws.Range("L5").Value = 1.66
$ws.Range("M5").Value = 1.88
$ws.Range("N5").Value = 0.5
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.08
$ws.Range("Q5").Value = 1.54
$ws.Range("R5").Value = 0.63
$ws.Range("S5").Value = 0.28
$ws.Range("T5").Value = 15.62
$ws.Range("U5").Value = 4.44
$ws.Range("V5").Value = 1.39
$ws.Range("W5").Value = 2.79
$ws.Range("X5").Value = 1.38
$ws.Range("Y5").Value = 0.47
$ws.Range("Z5").Value = 4.11
$ws.Range("AA5").Value = 1.23
$ws.Range("AB5").Value = 1.23
$ws.Range("AC5").Value = 1.41
$ws.Range("AD5").Value = 1.73
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 8.34
$ws.Range("AG5").Value = 0.65
$ws.Range("AH5").Value = 1.74

# Remove the last data row (row 6) - dataset trimmed to 1000 rows equivalent for this sheet
$ws.Rows("6:6").Delete()
